# Weekly refresh of the "Ají" (hortaliza) price-report subset.
# Each data row (2-17) gets updated to the newest weekly pull: date, variety,
# volume, min/max/avg price, commercialisation unit, $/Kg price and Kg-or-Units.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44553
$ws.Range("H2").Value = "Inferno"
$ws.Range("J2").Value = 35
$ws.Range("K2").Value = 45000
$ws.Range("L2").Value = 45000
$ws.Range("M2").Value = 45000
$ws.Range("P2").Value = 1800

$ws.Range("D3").Value = 44326
$ws.Range("K3").Value = 30000
$ws.Range("L3").Value = 30000
$ws.Range("M3").Value = 30000
$ws.Range("N3").Value = "`$/caja 25 kilos"
$ws.Range("P3").Value = 1200
$ws.Range("Q3").Value = 25

$ws.Range("D4").Value = 44474
$ws.Range("J4").Value = 18
$ws.Range("K4").Value = 100000
$ws.Range("L4").Value = 100000
$ws.Range("M4").Value = 100000
$ws.Range("P4").Value = 4000

$ws.Range("D5").Value = 44446
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 78000
$ws.Range("L5").Value = 78000
$ws.Range("M5").Value = 78000
$ws.Range("P5").Value = 3120

$ws.Range("D6").Value = 44446
$ws.Range("H6").Value = "Inferno"
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = 80000
$ws.Range("L6").Value = 80000
$ws.Range("M6").Value = 80000
$ws.Range("N6").Value = "`$/caja 15 kilos"
$ws.Range("P6").Value = 5333
$ws.Range("Q6").Value = 15

$ws.Range("D7").Value = 44425
$ws.Range("J7").Value = 15
$ws.Range("K7").Value = 75000
$ws.Range("L7").Value = 75000
$ws.Range("M7").Value = 75000
$ws.Range("P7").Value = 3000

$ws.Range("D8").Value = 44340
$ws.Range("H8").Value = "Americana (o)"
$ws.Range("J8").Value = 15

$ws.Range("D9").Value = 44449
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 80000
$ws.Range("L9").Value = 80000
$ws.Range("M9").Value = 80000
$ws.Range("P9").Value = 3200

$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 20
$ws.Range("K10").Value = 75000
$ws.Range("L10").Value = 75000
$ws.Range("M10").Value = 75000
$ws.Range("N10").Value = "`$/caja 15 kilos"
$ws.Range("P10").Value = 5000
$ws.Range("Q10").Value = 15

$ws.Range("D11").Value = 44421
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 15
$ws.Range("N11").Value = "`$/caja 25 kilos"
$ws.Range("P11").Value = 3000
$ws.Range("Q11").Value = 25

$ws.Range("D12").Value = 44319
$ws.Range("J12").Value = 20
$ws.Range("K12").Value = 30000
$ws.Range("L12").Value = 30000
$ws.Range("M12").Value = 30000
$ws.Range("P12").Value = 1200

$ws.Range("D13").Value = 44544
$ws.Range("J13").Value = 12
$ws.Range("K13").Value = 35000
$ws.Range("L13").Value = 35000
$ws.Range("M13").Value = 35000
$ws.Range("P13").Value = 1400

$ws.Range("D14").Value = 44343
$ws.Range("J14").Value = 20
$ws.Range("K14").Value = 36000
$ws.Range("L14").Value = 36000
$ws.Range("M14").Value = 36000
$ws.Range("P14").Value = 1440

$ws.Range("D15").Value = 44460
$ws.Range("J15").Value = 30
$ws.Range("K15").Value = 95000
$ws.Range("L15").Value = 95000
$ws.Range("M15").Value = 95000
$ws.Range("P15").Value = 3800

$ws.Range("D16").Value = 44221
$ws.Range("H16").Value = "Americana (o)"
$ws.Range("J16").Value = 22
$ws.Range("K16").Value = 24000
$ws.Range("L16").Value = 25000
$ws.Range("M16").Value = 24545
$ws.Range("N16").Value = "`$/caja 25 kilos"
$ws.Range("P16").Value = 982
$ws.Range("Q16").Value = 25

$ws.Range("D17").Value = 44193
$ws.Range("K17").Value = 46000
$ws.Range("L17").Value = 46000
$ws.Range("M17").Value = 46000
$ws.Range("N17").Value = "`$/caja 25 kilos"
$ws.Range("P17").Value = 3067
$ws.Range("Q17").Value = 15
